$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.286.81'
$ws.Range('E2').Value = '  +1.86%  '
$ws.Range('D3').Value = '1.645.84'
$ws.Range('E3').Value = '  +0.49%  '
$ws.Range('E4').Value = '  +0.03%  '
$ws.Range('D5').Value = '''217.54'
$ws.Range('E5').Value = '  +0.37%  '
$ws.Range('E6').Value = '  +1.90%  '
$ws.Range('E7').Value = '  +0.06%  '
$ws.Range('E8').Value = '  +0.41%  '
$ws.Range('E9').Value = '  +1.20%  '
$ws.Range('E10').Value = '  +1.11%  '
$ws.Range('D11').Value = '''0.0851'
$ws.Range('E11').Value = '  +0.71%  '
$ws.Range('D12').Value = '1.878.00'
$ws.Range('E12').Value = '  +0.59%  '
$ws.Range('D13').Value = '1.637.83'
$ws.Range('E13').Value = '  -0.12%  '
$ws.Range('D14').Value = '''4.15'
$ws.Range('E14').Value = '  +1.07%  '
$ws.Range('D15').Value = '''0.546'
$ws.Range('E15').Value = '  +3.30%  '
$ws.Range('D16').Value = '''67.19'
$ws.Range('E16').Value = '  +1.20%  '
$ws.Range('D17').Value = '27.288.66'
$ws.Range('E17').Value = '  +1.84%  '
$ws.Range('D18').Value = '0.0₃0741'
$ws.Range('E18').Value = '  +1.80%  '
$ws.Range('D19').Value = '''220.32'
$ws.Range('E19').Value = '  +0.69%  '
$ws.Range('E20').Value = '  -0.02%  '
$ws.Range('D21').Value = '''6.98'
$ws.Range('E21').Value = '  +4.48%  '
$ws.Range('E22').Value = '  +3.76%  '
$ws.Range('D23').Value = '''4.42'
$ws.Range('E23').Value = '  +0.74%  '
$ws.Range('D24').Value = '''9.16'
$ws.Range('E24').Value = '  +0.35%  '
$ws.Range('D25').Value = '''148.34'
$ws.Range('E25').Value = '  +0.80%  '
$ws.Range('B26').Value = 'BinanceUSD'
$ws.Range('C26').Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range('D26').Value = '''1.00'
$ws.Range('E26').Value = '  +0.03%  '
$ws.Range('B27').Value = 'Cosmos'
$ws.Range('C27').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D27').Value = '''7.49'
$ws.Range('E27').Value = '  +1.55%  '
$ws.Range('E28').Value = '  -0.04%  '
$ws.Range('D29').Value = '''15.77'
$ws.Range('E29').Value = '  +0.29%  '
$ws.Range('D30').Value = '''0.0512'
$ws.Range('E30').Value = '  +1.80%  '
$ws.Range('E31').Value = '  +0.65%  '
$ws.Range('E32').Value = '  +2.08%  '
$ws.Range('E33').Value = '  +0.74%  '
$ws.Range('E34').Value = '  +1.54%  '
$ws.Range('D35').Value = '1.303.37'
$ws.Range('E35').Value = '  +3.41%  '
$ws.Range('E36').Value = '  +1.34%  '
$ws.Range('D37').Value = '''0.0177'
$ws.Range('E37').Value = '  -0.02%  '
$ws.Range('E38').Value = '  +3.72%  '
$ws.Range('D39').Value = '''0.861'
$ws.Range('E39').Value = '  +3.36%  '
$ws.Range('E40').Value = '  +0.05%  '
$ws.Range('D41').Value = '''0.811'
$ws.Range('E41').Value = '  +0.54%  '
$ws.Range('D42').Value = '''2.23'
$ws.Range('E42').Value = '  +6.00%  '
$ws.Range('E43').Value = '  -2.45%  '
$ws.Range('D44').Value = '1.787.26'
$ws.Range('E44').Value = '  +0.55%  '
$ws.Range('D45').Value = '''62.20'
$ws.Range('E45').Value = '  +0.79%  '
$ws.Range('D46').Value = '''92.14'
$ws.Range('E46').Value = '  +0.51%  '
$ws.Range('D47').Value = '''1.60'
$ws.Range('E47').Value = '  +2.34%  '
$ws.Range('D48').Value = '0.0₆0107'
$ws.Range('E48').Value = '  +1.24%  '
$ws.Range('D49').Value = '''0.0513'
$ws.Range('D50').Value = '''7.69'
$ws.Range('E50').Value = '  +0.99%  '
$ws.Range('D51').Value = '''0.0968'
$ws.Range('E51').Value = '  +0.51%  '
